$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.445.70'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '2.405.96'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '506.89'
$ws.Range('E5').Value = '  -3.43%  '
$ws.Range('D6').Value = '133.37'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').Value = '0.994'
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').Value = '2.442.52'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '0.0979'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = '4.60'
$ws.Range('E13').Value = '  -7.55%  '
$ws.Range('D14').Value = '2.843.35'
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').Value = '57.273.62'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '21.93'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '2.405.89'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('D19').Value = '10.29'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').Value = '4.11'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('D21').Value = '313.86'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '6.41'
$ws.Range('E22').Value = '  +5.23%  '
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').Value = '65.14'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').Value = '0.992'
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').Value = '2.527.51'
$ws.Range('E27').Value = '  -6.60%  '
$ws.Range('D28').Value = '0.380'
$ws.Range('E28').Value = '  -5.84%  '
$ws.Range('E29').Value = '  -2.22%  '
$ws.Range('D30').Value = '7.59'
$ws.Range('E30').Value = '  +4.84%  '
$ws.Range('D31').Value = '173.48'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').Value = '0.0₃0732'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').Value = '6.18'
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('E39').Value = '  +3.15%  '
$ws.Range('D40').Value = '3.82'
$ws.Range('E40').Value = '  +0.41%  '
$ws.Range('D41').Value = '36.56'
$ws.Range('D42').Value = '0.814'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('D44').Value = '134.82'
$ws.Range('E44').Value = '  +10.22%  '
$ws.Range('E45').Value = '  +3.97%  '
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('D47').Value = '256.27'
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('E48').Value = '  -2.97%  '
$ws.Range('D49').Value = '0.0916'
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').Value = '0.0214'
$ws.Range('E51').Value = '  +0.63%  '
